$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista de Itens de Trabalho")

# --- Row 30: new work item "Corrigir todos os documentos para o início de GPSII" ---
$ws.Range("A30").Value = "Corrigir todos os documentos para o início de GPSII"
$ws.Range("B30").Value = "Alta"
$ws.Range("C30").Value = 100

# D30 needs to stay a literal text string "0.5%" while keeping its existing
# percentage number format (style index 7). A direct .Value assignment gets
# auto-coerced into a real percentage number by Excel, so instead we build it
# as a text formula and then paste the computed value back in-place, which
# preserves the original cell style/number format instead of re-deriving one.
$ws.Range("D30").Formula = "=""0.5%"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("F30").Value = "Valter "
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = 0.5

# I30 holds the same GitHub repo link used throughout column I.
$ws.Range("I30").Value = "https://github.com/valtervasconcelos/Projeto_de_Software_1"
$ws.Hyperlinks.Add($ws.Range("I30"), "https://github.com/valtervasconcelos/Projeto_de_Software_1") | Out-Null
# Hyperlinks.Add forces Excel's built-in blue "Hyperlink" style; restore the
# custom color/underline already used by the other hyperlink cells in I.
$ws.Range("I30").Font.Color = $ws.Range("I11").Font.Color
$ws.Range("I30").Font.Underline = $ws.Range("I11").Font.Underline

# --- Selection recorded in the saved view state ---
$ws.Range("H30").Select()
